$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "reset options?" prompt text (row 25, column B)
$ws.Range("B25").Value = "Reset options and control bindings?"

# 2. Update "reset" button text (row 26, column B)
$ws.Range("B26").Value = "Reset"

# 3. Insert a new row after "recenter pose" (row 68) for the new
#    "recenter pose button" VR control-hint string, then populate it.
$ws.Rows(69).Insert()
$ws.Range("A69").Value = "recenter pose button"
$ws.Range("B69").Value = "Recenter VR pose [{{RecenterVRPose}}]"

# 4. Update the saved sheet selection/scroll position.
[void]$ws.Range("B27").Select()
